$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the "Actual Time" values that were just measured for the
#     bomb-placement / collision-check / brick-destruction tasks ---
$ws.Range("E9").Value  = 1
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 0.5

# --- Make room for a new task row above the "Total" row: push the
#     Total row from 18 down to 19, carrying its formatting with it,
#     then rebuild its content/formula in place ---
$ws.Range("C18:E18").Copy()
$ws.Range("C19:E19").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C19").Value = "Total"
$ws.Range("D19").Formula = "=SUM(D5:D18)"
$ws.Range("E19").Formula = "=SUM(E5:E17)"
$ws.Rows("19:19").RowHeight = 18.75

# --- Turn (the now-vacated) row 18 into the new optional "Bomb & Brick
#     FX" task row, copying the formatting of the row above it so it
#     matches the rest of the table (borders, centering, wrap text) ---
$ws.Range("C17:E17").Copy()
$ws.Range("C18:E18").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C18").Value = "Bomb & Brick FX(Optional)"
$ws.Range("D18").Value = 2
$ws.Range("E18").ClearContents()
$ws.Rows("18:18").RowHeight = 30

# Restore the selection to what was recorded after the edit.
$ws.Range("J17").Select() | Out-Null
